# #355 basic attribute generation works
# Adds a new "Attribute" sheet (after the existing "StructureOrder" sheet)
# that documents a struct called "Attribute" with two string attributes,
# mirroring the layout used by the other schema sheets in this workbook.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet right after the current last sheet so it lands
# at the end of the tab strip (becomes the new active/selected sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Attribute"

# Header row
$ws.Range("A1").Value = "class"
$ws.Range("B1").Value = "name"
$ws.Range("A1:B1").Font.Bold = $true

# Data rows
$ws.Range("A2").Value = "struct"
$ws.Range("B2").Value = "Attribute"

$ws.Range("A3").Value = "attribute"
$ws.Range("B3").Value = "stringAttr1"

$ws.Range("A4").Value = "attribute"
$ws.Range("B4").Value = "stringAttr2"

# Column widths to match the other schema sheets in the workbook
# (the host quantizes ColumnWidth onto a 1/6-character grid, so these
# values are chosen to land on the closest achievable gridpoints to the
# authored widths of 16.84375 / 15.3828125 / 12.53515625).
$ws.Columns.Item(2).ColumnWidth = 16.0
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 11.7

# Leave the selection on A8, matching the authored sheet.
$ws.Range("A8").Select()
